$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header fields ---
$ws.Range("C2").Value = "Hartmut"

# B3 holds a purely numeric-looking card number that must stay TEXT
# (matching the source data, which stores it as an inline string).
# Assigning a digit-only string straight to .Value would make Excel
# auto-coerce it to a number, so force a text result via a formula
# and then flatten the formula down to a static value.
$ws.Range("B3").Formula = '=T("2570314725427075")'
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 23.01.2024"

# --- Row 6 ---
$ws.Range("B6").Value = "25.01."
$ws.Range("C6").Value = "26.01."
$ws.Range("D6").Value = "PAYPAL NZFTKS"
$ws.Range("E6").Value = "72,08-"

# --- Row 7 ---
$ws.Range("B7").Value = "27.01."
$ws.Range("C7").Value = "28.01."
$ws.Range("D7").Value = "AMAZON.DE MKTPLC EU VNWRKO"
$ws.Range("E7").Value = "20,12-"

# --- Row 8 ---
$ws.Range("B8").Value = "28.01."
$ws.Range("C8").Value = "29.01."
$ws.Range("D8").Value = "RECHNUNG VODAFONE GMBH 22884255"
$ws.Range("E8").Value = "40,44-"

# --- Row 9: previously blank, now a new transaction row ---
# Copy formatting from E8 (style 17, right aligned) onto E9 so the
# cell style matches the target (was style 13, now style 17).
$ws.Range("E8").Copy() | Out-Null
$ws.Range("E9").PasteSpecial(-4122) | Out-Null
$ws.Range("B9").Value = "01.02."
$ws.Range("C9").Value = "02.02."
$ws.Range("D9").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E9").Value = "24,58-"

# --- Row 10: previously blank, now a new transaction row ---
$ws.Range("E8").Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Value = "02.02."
$ws.Range("C10").Value = "03.02."
$ws.Range("D10").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E10").Value = "53,40-"

$excel.CutCopyMode = $false

# --- Closing balance / footer ---
$ws.Range("D12").Value = "KONTOSTAND AM 06.02.2024"
$ws.Range("E12").Value = "210,62-"
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 11.02.2024"
